# Updates crypto price/volume figures (and one Maker/Bittensor row swap)
# as captured by the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.886.71'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '3.727.91'
$ws.Range("E3").Value = '  +18.58%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''614.07'
$ws.Range("E5").Value = '  +6.19%  '
$ws.Range("D6").Value = '''176.78'
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("D7").Value = '3.725.65'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +3.30%  '
$ws.Range("E10").Value = '  +9.16%  '
$ws.Range("D11").Value = '''6.38'
$ws.Range("E11").Value = '  -2.05%  '
$ws.Range("D12").Value = '''0.499'
$ws.Range("E12").Value = '  +5.97%  '
$ws.Range("D13").Value = '''40.68'
$ws.Range("E13").Value = '  +10.11%  '
$ws.Range("E14").Value = '  +5.29%  '
$ws.Range("D15").Value = '4.353.94'
$ws.Range("E15").Value = '  +18.56%  '
$ws.Range("D16").Value = '3.729.08'
$ws.Range("E16").Value = '  +18.55%  '
$ws.Range("D17").Value = '69.962.59'
$ws.Range("E17").Value = '  +1.77%  '
$ws.Range("D18").Value = '''0.123'
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  +6.16%  '
$ws.Range("D20").Value = '''515.18'
$ws.Range("E20").Value = '  +5.22%  '
$ws.Range("D21").Value = '''16.68'
$ws.Range("E21").Value = '  +1.65%  '
$ws.Range("D22").Value = '''9.42'
$ws.Range("E22").Value = '  +20.97%  '
$ws.Range("D23").Value = '''0.727'
$ws.Range("E23").Value = '  +4.35%  '
$ws.Range("D24").Value = '''88.16'
$ws.Range("E24").Value = '  +4.77%  '
$ws.Range("E25").Value = '  +5.97%  '
$ws.Range("E26").Value = '  +4.53%  '
$ws.Range("D27").Value = '''10.98'
$ws.Range("E27").Value = '  +3.84%  '
$ws.Range("D29").Value = '''0.0000126'
$ws.Range("E29").Value = '  +32.42%  '
$ws.Range("E30").Value = '  +6.42%  '
$ws.Range("E31").Value = '  +7.83%  '
$ws.Range("E32").Value = '  -3.35%  '
$ws.Range("D33").Value = '''31.33'
$ws.Range("E33").Value = '  +11.25%  '
$ws.Range("E34").Value = '  +2.58%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  +7.72%  '
$ws.Range("D37").Value = '''1.04'
$ws.Range("E37").Value = '  +7.81%  '
$ws.Range("D38").Value = '''0.340'
$ws.Range("E38").Value = '  +5.35%  '
$ws.Range("E39").Value = '  +6.87%  '
$ws.Range("E40").Value = '  +5.61%  '
$ws.Range("D41").Value = '''51.24'
$ws.Range("E41").Value = '  +4.04%  '
$ws.Range("D42").Value = '''8.81'
$ws.Range("E42").Value = '  +5.38%  '
$ws.Range("D43").Value = '''44.22'
$ws.Range("E43").Value = '  -7.71%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.076.67'
$ws.Range("E44").Value = '  +9.33%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '''418.99'
$ws.Range("E45").Value = '  +5.25%  '
$ws.Range("D46").Value = '''2.71'
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("D47").Value = '''0.0364'
$ws.Range("E47").Value = '  +4.84%  '
$ws.Range("D48").Value = '''27.85'
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("E49").Value = '  +6.25%  '
$ws.Range("D50").Value = '''135.93'
$ws.Range("E50").Value = '  +0.50%  '
